$wb = $excel.ActiveWorkbook

# 1. Rename the "Data" sheet to "Data  for Class 4"
$wsData = $wb.Worksheets.Item("Data")
$wsData.Name = "Data  for Class 4"

# 2. On the "Class-4" sheet, remove the duplicated "Data" table (old columns
#    A:J) so only the small numbered list that used to live in columns L:M
#    remains - it shifts left into columns B:C.
$wsClass4 = $wb.Worksheets.Item("Class-4")
$wsClass4.Range("A1:J1").EntireColumn.Delete()

# 3. Update the visible selection on Class-4 to E11 without stealing the
#    "active sheet" state from the Data sheet (match original tab order).
$wsClass4.Select()
$wsClass4.Range("E11").Select()
$wsData.Select()
